# Adds Endotherm to KillRatePreyMassExponent in RevisedPredation
# New row 22 mirrors the existing "formula" rows (e.g. row 10/11/13/14):
#   D = internal variable name
#   E = formula path (RevisedPredation.KillRatePreyMassExponent)
#   F = Endotherm-specific formula path (RevisedPredation.Endotherm.KillRatePreyMassExponent)
#   H = numeric value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = "_KillRatePreyMassExponent"
$ws.Range("E22").Value = "Predation.RevisedPredation.KillRatePreyMassExponent"
$ws.Range("F22").Value = "Predation.RevisedPredation.Endotherm.KillRatePreyMassExponent"
$ws.Range("H22").Value = -0.08832

# Match the formatting already used on the other "Value"/formula columns (E/F) -
# red font, same as E10:F14.
$ws.Range("E22:F22").Font.Color = 255

# Move the selection to the newly-added row, as the author did.
$ws.Range("F22").Select()
